# Add season record (Wins/Losses/Ties) columns to the roster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the other header cells (bold, bordered,
# centered/top-aligned) by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's season record for every player row.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 68
    $ws.Cells.Item($row, 31).Value = 94
    $ws.Cells.Item($row, 32).Value = 0
}
